$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the footnote texts in row 8 (source attribution changed) ---
# Order matches the author's edit: B8, C8, A8 (per shared-strings layout)
$ws.Range("B8").Value = "*по данным МЦР КР"
$ws.Range("C8").Value = "*according to the MDD KR"
$ws.Range("A8").Value = "*КР СӨМ маалыматтары  боюнча"

# --- Add a new year column O (2023 data), cloning per-row formatting from column N ---
$ws.Range("N3").Copy()
$ws.Range("O3").PasteSpecial(-4122)

$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial(-4122)

$ws.Range("N5").Copy()
$ws.Range("O5").PasteSpecial(-4122)

$ws.Range("N6").Copy()
$ws.Range("O6").PasteSpecial(-4122)

$ws.Range("N7").Copy()
$ws.Range("O7").PasteSpecial(-4122)

$excel.CutCopyMode = $false

$ws.Range("O4").Value = 2023
$ws.Range("O5").Value = 99
$ws.Range("O6").Value = 98.9
$ws.Range("O7").Value = 98.8

# --- Widen columns A:C from 35.71 to 38 characters ---
$ws.Columns("A:C").ColumnWidth = 37.166666666666664
